$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the envelope-related column headers to use underscores instead of spaces
$ws.Range("J1").Value = "Window_type"
$ws.Range("K1").Value = "Window_area"
$ws.Range("L1").Value = "Door_type"
$ws.Range("M1").Value = "Door_area"
$ws.Range("N1").Value = "Net_wall_area"
$ws.Range("O1").Value = "Ceiling/attic_area"
$ws.Range("P1").Value = "Basment_type"
$ws.Range("Q1").Value = "Basement_area"

# Update the rendering/view setup: zoom level and the active selection
$ws.Application.ActiveWindow.Zoom = 125
$ws.Range("H4").Select()
